$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("立讯精密", "立讯精密", "首开股份"),
    @("山子高科", "山子高科", "山河智能"),
    @("工业富联", "和而泰", "领益智造"),
    @("和而泰", "上海建工", "长盈精密"),
    @("上海建工", "工业富联", "上海建工"),
    @("领益智造", "首开股份", "卧龙电驱"),
    @("首开股份", "领益智造", "先导智能"),
    @("东华软件", "东华软件", "利欧股份"),
    @("凯美特气", "凯美特气", "张江高科"),
    @("卧龙电驱", "卧龙电驱", "工业富联"),
    @("山河智能", "中科曙光", "三花智控"),
    @("福龙马", "福龙马", "国盛金控"),
    @("利欧股份", "歌尔股份", "山子高科"),
    @("三花智控", "山河智能", "华胜天成"),
    @("歌尔股份", "利欧股份", "立讯精密"),
    @("初灵信息", "三花智控", "欧菲光"),
    @("中科曙光", "盈趣科技", "指南针"),
    @("长盈精密", "万马股份", "吉视传媒"),
    @("张江高科", "长盈精密", "西部黄金"),
    @("中电鑫龙", "张江高科", "晓程科技")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
